# Add a default (primary) header and footer to the document's only section.
# Touching the HeaderFooter's Range (rather than assigning .Text, which would
# cause Word to materialize all six header/footer slots - primary, even page
# and first page for both header and footer) keeps this to exactly the
# "default" slot, matching a simple "insert blank header/footer" edit.
$d = $word.ActiveDocument
$sec = $d.Sections.First

$header = $sec.Headers.Item(1)
$header.Range.Style = "Header"

$footer = $sec.Footers.Item(1)
$footer.Range.Style = "Footer"
